$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The C3DC queries referenced the old "id" / "participant.id" / "study.id"
# column names. The source dataframes were renamed to use explicit
# "study_id" / "participant_id" columns, so every JOIN clause across the
# embedded SQL queries needs updating to match.
$replacements = @(
    @('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"'),
    @('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"'),
    @('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"'),
    @('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"'),
    @('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"'),
    @('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
)

# Every query cell on the sheet (StatQuery in C2, TabQuery in B2:B7) contains
# the same set of LEFT JOIN clauses that need the substitutions above.
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($cellAddr in $cells) {
    $range = $ws.Range($cellAddr)
    $text = $range.Value()
    foreach ($pair in $replacements) {
        $text = $text -replace [regex]::Escape($pair[0]), $pair[1]
    }
    $range.Value = $text
}

# Widen column C (StatQuery) from its old auto-fit width to a fixed 69
# characters, matching the new, longer query text.
$ws.Columns.Item(3).ColumnWidth = 68.16666666666667
